$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "contextId"
$ws.Range("B1").Value = "cueId"
$ws.Range("C1").Value = "sick"
$ws.Range("D1").Value = "corrAns"

# --- Row 2: "Training" section header (merged A2:D2) ---
$ws.Range("A2").Value = "Training"

# --- Row 3 ---
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = "Yes"
$ws.Range("D3").Value = "left"

# --- Row 4 ---
$ws.Range("A4").Value = 0
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "No"
$ws.Range("D4").Value = "right"

# --- Row 5 ---
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = "Yes"
$ws.Range("D5").Value = "left"

# --- Row 6 ---
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = "No"
$ws.Range("D6").Value = "right"

# --- Row 7: " Test" section header (merged A7:D7) ---
$ws.Range("A7").Value = " Test"

# --- Row 8 ---
$ws.Range("A8").Value = 0
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = "Yes"
$ws.Range("D8").Value = "left"

# --- Row 9 ---
$ws.Range("A9").Value = 2
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = "Yes"
$ws.Range("D9").Value = "left"

# --- Row 10 ---
$ws.Range("A10").Value = 0
$ws.Range("B10").Value = 2

# --- Row 11 ---
$ws.Range("A11").Value = 2
$ws.Range("B11").Value = 2

# --- Update the active selection to C10 ---
$ws.Range("C10").Select()
